$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1633.125
$ws.Range("I15").Value = 1633.125
$ws.Range("K15").Value = 4899.375
$ws.Range("M15").Value = -4730.375
$ws.Range("H17").Value = 611.2456
$ws.Range("J17").Value = 614.1429000000001
$ws.Range("L17").Value = 1842.4287
$ws.Range("N17").Value = -2178.4287
$ws.Range("H18").Value = 3592.2856
$ws.Range("I18").Value = 1724.625
$ws.Range("J18").Value = 6082.5
$ws.Range("K18").Value = 1724.625
$ws.Range("L18").Value = 6082.5
$ws.Range("M18").Value = -1440.625
$ws.Range("N18").Value = -6650.5
$ws.Range("H28").Value = 837.34784
$ws.Range("I28").Value = 466.47058
$ws.Range("J28").Value = 1888.1666
$ws.Range("K28").Value = 466.47058
$ws.Range("L28").Value = 1888.1666
$ws.Range("M28").Value = 18.52942000000002
$ws.Range("N28").Value = -2858.1666
$ws.Range("H33").Value = 1226.7273
$ws.Range("I33").Value = 1263.15
$ws.Range("K33").Value = 1263.15
$ws.Range("M33").Value = -1034.15
$ws.Range("H88").Value = 700
$ws.Range("I88").Value = 700
$ws.Range("K88").Value = 700
$ws.Range("M88").Value = -294
$ws.Range("H91").Value = 700
$ws.Range("I91").Value = 700
$ws.Range("K91").Value = 700
$ws.Range("M91").Value = 704
$ws.Range("H98").Value = 3428.7026
$ws.Range("I98").Value = 1286.12
$ws.Range("J98").Value = 7892.4165
$ws.Range("K98").Value = 1286.12
$ws.Range("L98").Value = 7892.4165
$ws.Range("M98").Value = 211.8800000000001
$ws.Range("N98").Value = -10888.4165
$ws.Range("H101").Value = 2499.1667
$ws.Range("J101").Value = 2749
$ws.Range("L101").Value = 8247
$ws.Range("N101").Value = -11491
$ws.Range("H116").Value = 2920.2666
$ws.Range("I116").Value = 2770.5
$ws.Range("J116").Value = 3219.8
$ws.Range("K116").Value = 2770.5
$ws.Range("L116").Value = 3219.8
$ws.Range("M116").Value = 671.5
$ws.Range("N116").Value = -10103.8
$ws.Range("H122").Value = 3428.7026
$ws.Range("I122").Value = 1286.12
$ws.Range("J122").Value = 7892.4165
$ws.Range("K122").Value = 3858.36
$ws.Range("L122").Value = 23677.2495
$ws.Range("M122").Value = -1408.36
$ws.Range("N122").Value = -28577.2495
$ws.Range("H132").Value = 12990141
$ws.Range("I132").Value = 16396268
$ws.Range("J132").Value = 4281.4375
$ws.Range("K132").Value = 49188804
$ws.Range("L132").Value = 12844.3125
$ws.Range("M132").Value = -49186274
$ws.Range("N132").Value = -17904.3125
$ws.Range("H138").Value = 1790.8572
$ws.Range("I138").Value = 1083.3684
$ws.Range("J138").Value = 8512
$ws.Range("K138").Value = 3250.1052
$ws.Range("L138").Value = 25536
$ws.Range("M138").Value = 1889.8948
$ws.Range("N138").Value = -35816
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2088.9565
$ws.Range("I2").Value = 2168.8572
$ws.Range("K2").Value = 2168.8572
$ws.Range("M2").Value = -2055.8572
$ws.Range("H31").Value = 11009.5
$ws.Range("J31").Value = 20519
$ws.Range("L31").Value = 20519
$ws.Range("N31").Value = -21107
$ws.Range("H32").Value = 3154.7292
$ws.Range("I32").Value = 3136.4468
$ws.Range("K32").Value = 3136.4468
$ws.Range("M32").Value = -2849.4468
$ws.Range("H74").Value = 868.7742
$ws.Range("I74").Value = 622.2273
$ws.Range("K74").Value = 622.2273
$ws.Range("M74").Value = 251.7727
$ws.Range("H77").Value = 868.7742
$ws.Range("I77").Value = 622.2273
$ws.Range("K77").Value = 3111.1365
$ws.Range("M77").Value = 1256.8635
$ws.Range("H93").Value = 66882.60000000001
$ws.Range("J93").Value = 66882.60000000001
$ws.Range("L93").Value = 66882.60000000001
$ws.Range("N93").Value = -71874.60000000001
$ws.Range("H116").Value = 2088.9565
$ws.Range("I116").Value = 2168.8572
$ws.Range("K116").Value = 2168.8572
$ws.Range("M116").Value = 125.1428000000001
$ws.Range("H122").Value = 5139.2583
$ws.Range("I122").Value = 4916.304
$ws.Range("K122").Value = 14748.912
$ws.Range("M122").Value = -12298.912
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2088.9565
$ws.Range("I3").Value = 2168.8572
$ws.Range("K3").Value = 2168.8572
$ws.Range("M3").Value = -2054.8572
$ws.Range("H82").Value = 17744.084
$ws.Range("I82").Value = 4756.222
$ws.Range("J82").Value = 56707.668
$ws.Range("K82").Value = 4756.222
$ws.Range("L82").Value = 56707.668
$ws.Range("M82").Value = -4373.222
$ws.Range("N82").Value = -57473.668
$ws.Range("H85").Value = 17744.084
$ws.Range("I85").Value = 4756.222
$ws.Range("J85").Value = 56707.668
$ws.Range("K85").Value = 4756.222
$ws.Range("L85").Value = 56707.668
$ws.Range("M85").Value = -3430.222
$ws.Range("N85").Value = -59359.668
$ws.Range("H93").Value = 66666.336
$ws.Range("I93").Value = 49999
$ws.Range("K93").Value = 49999
$ws.Range("M93").Value = -48127
$ws.Range("H94").Value = 2499.8333
$ws.Range("I94").Value = 3800
$ws.Range("K94").Value = 3800
$ws.Range("M94").Value = -3349
$ws.Range("H99").Value = 2409.3635
$ws.Range("I99").Value = 2450.2
$ws.Range("K99").Value = 2450.2
$ws.Range("M99").Value = -952.1999999999998
$ws.Range("H105").Value = 3410.8125
$ws.Range("I105").Value = 3454
$ws.Range("J105").Value = 3377.2222
$ws.Range("K105").Value = 3454
$ws.Range("L105").Value = 3377.2222
$ws.Range("M105").Value = -1707
$ws.Range("N105").Value = -6871.2222
$ws.Range("H134").Value = 2962.6191
$ws.Range("I134").Value = 2460.5
$ws.Range("K134").Value = 7381.5
$ws.Range("M134").Value = -4846.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 430.67648
$ws.Range("I22").Value = 389.17392
$ws.Range("K22").Value = 389.17392
$ws.Range("M22").Value = -39.17392000000001
$ws.Range("H105").Value = 2862
$ws.Range("I105").Value = 2862
$ws.Range("K105").Value = 2862
$ws.Range("M105").Value = -1115
$ws.Range("H132").Value = 1580.1666
$ws.Range("I132").Value = 1758.2
$ws.Range("K132").Value = 5274.6
$ws.Range("M132").Value = -2744.6
$ws.Range("H141").Value = 79513.39999999999
$ws.Range("I141").Value = 71715.2
$ws.Range("K141").Value = 71715.2
$ws.Range("M141").Value = -66535.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1109.2
$ws.Range("I50").Value = 676.8889
$ws.Range("K50").Value = 2030.6667
$ws.Range("M50").Value = -1549.6667
$ws.Range("H53").Value = 1109.2
$ws.Range("I53").Value = 676.8889
$ws.Range("K53").Value = 2030.6667
$ws.Range("M53").Value = -1549.6667
$ws.Range("H122").Value = 756.2778
$ws.Range("J122").Value = 903.8333
$ws.Range("L122").Value = 8134.4997
$ws.Range("N122").Value = -13034.4997
$ws.Range("H132").Value = 2138.8147
$ws.Range("I132").Value = 1373.5
$ws.Range("J132").Value = 2271.913
$ws.Range("K132").Value = 12361.5
$ws.Range("L132").Value = 20447.217
$ws.Range("M132").Value = -9831.5
$ws.Range("N132").Value = -25507.217
$ws.Range("H140").Value = 4538
$ws.Range("I140").Value = 3715.1
$ws.Range("K140").Value = 11145.3
$ws.Range("M140").Value = -5965.299999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6928.909
$ws.Range("I80").Value = 7165.643
$ws.Range("K80").Value = 7165.643
$ws.Range("M80").Value = -6167.643
$ws.Range("H83").Value = 6928.909
$ws.Range("I83").Value = 7165.643
$ws.Range("K83").Value = 35828.215
$ws.Range("M83").Value = -30836.215
$ws.Range("H102").Value = 74381.86
$ws.Range("I102").Value = 2993.6667
$ws.Range("J102").Value = 127923
$ws.Range("K102").Value = 2993.6667
$ws.Range("L102").Value = 127923
$ws.Range("M102").Value = -1371.6667
$ws.Range("N102").Value = -131167
$ws.Range("H123").Value = 38722.223
$ws.Range("J123").Value = 38722.223
$ws.Range("L123").Value = 38722.223
$ws.Range("N123").Value = -43622.223
$ws.Range("H126").Value = 5227.9575
$ws.Range("I126").Value = 4489.2354
$ws.Range("K126").Value = 13467.7062
$ws.Range("M126").Value = -10997.7062
$ws.Range("H132").Value = 6787.289
$ws.Range("I132").Value = 7528.9473
$ws.Range("K132").Value = 22586.8419
$ws.Range("M132").Value = -20056.8419
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5862.4375
$ws.Range("I40").Value = 5138.385
$ws.Range("J40").Value = 9000
$ws.Range("K40").Value = 5138.385
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = -5002.385
$ws.Range("N40").Value = -9272
$ws.Range("H122").Value = 15268.477
$ws.Range("I122").Value = 14415
$ws.Range("J122").Value = 17673.727
$ws.Range("K122").Value = 43245
$ws.Range("L122").Value = 53021.181
$ws.Range("M122").Value = -40795
$ws.Range("N122").Value = -57921.181
$ws.Range("H132").Value = 4298.7896
$ws.Range("I132").Value = 2973.3333
$ws.Range("K132").Value = 8919.999899999999
$ws.Range("M132").Value = -6389.999899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 52663.24
$ws.Range("I81").Value = 102159.6
$ws.Range("J81").Value = 7666.5454
$ws.Range("K81").Value = 204319.2
$ws.Range("L81").Value = 15333.0908
$ws.Range("M81").Value = -203258.2
$ws.Range("N81").Value = -17455.0908
$ws.Range("H84").Value = 52663.24
$ws.Range("I84").Value = 102159.6
$ws.Range("J84").Value = 7666.5454
$ws.Range("K84").Value = 1021596
$ws.Range("L84").Value = 76665.454
$ws.Range("M84").Value = -1016292
$ws.Range("N84").Value = -87273.454
$ws.Range("H132").Value = 1828.5714
$ws.Range("J132").Value = 1600
$ws.Range("L132").Value = 4800
$ws.Range("N132").Value = -9860

Write-Host "Edit complete"